# Generate Report for Handoff
#
# The localization-status report is regenerated: the "Status" for the
# zh-cn / de-de targets moves from "Handed back: in sync with en-US" to
# "Ready for handoff", the handoff timestamps advance, and the Status
# columns on the per-language sheets are narrowed to match the new
# (shorter) status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
# Overview!E2 (zh-cn status), Overview!F2 (de-de status)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# zh-cn!C2 and de-de!C2 (Status column on each language sheet)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2 "Latest Handoff
# Datetime" both advance from 2016-09-06 01:05:20 -> 2016-09-06 01:06:03
$wsOverview.Range("G2").Value = "2016-09-06 01:06:03"
$wsDeDe.Range("H2").Value = "2016-09-06 01:06:03"

# zh-cn!H2 "Latest Handoff Datetime" advances from
# 2016-09-06 01:05:15 -> 2016-09-06 01:05:57
$wsZhCn.Range("H2").Value = "2016-09-06 01:05:57"

# --- Column widths ---
# The Status columns are narrowed (report regenerated with the shorter
# "Ready for handoff" text instead of "Handed back: in sync with en-US").
# Overview columns E (zh-cn status) and F (de-de status); zh-cn/de-de
# column C (Status).
$wsOverview.Range("E1:F1").ColumnWidth = 16.3
$wsZhCn.Range("C1").ColumnWidth = 16.3
$wsDeDe.Range("C1").ColumnWidth = 16.3
